$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Changes 1-4: pure run-merges (identical text, just re-typed as one run).
# Replacing the matched text with itself causes the engine to coalesce the
# previously-split (but identically formatted) runs into a single run.
# ---------------------------------------------------------------------------

# Luke 4:18-19 scripture quote
$d.Content.Find.Execute("anointed me to proclaim", $true, $false, $false, $false, $false, $true, 1, $false, "anointed me to proclaim", 2) | Out-Null

# "DISCUSSION " + "QUESTIONS"
$d.Content.Find.Execute("DISCUSSION QUESTIONS", $true, $false, $false, $false, $false, $true, 1, $false, "DISCUSSION QUESTIONS", 2) | Out-Null

# "found som" + "e purpose"
$d.Content.Find.Execute("found some purpose", $true, $false, $false, $false, $false, $true, 1, $false, "found some purpose", 2) | Out-Null

# "you can" + " move forward"
$d.Content.Find.Execute("you can move forward", $true, $false, $false, $false, $false, $true, 1, $false, "you can move forward", 2) | Out-Null

# ---------------------------------------------------------------------------
# fpStudents / Christmas paragraph: extend the break to two weeks, move the
# "back on" date from Jan. 1st to Jan. 8th. Do this BEFORE the Fusion
# paragraph edit below so the (still valid, original-document) character
# offsets used here are not disturbed by the Fusion paragraph's own edits,
# which sit earlier in the story and would otherwise shift everything after
# them.
# ---------------------------------------------------------------------------

# "st" (superscript) -> "th" (keep superscript formatting, only change text)
$fpSuperscript = $d.Range(1716, 1718)
$fpSuperscript.Text = "th"

# Main sentence text, up to (but excluding) the old "st"
$fpMain = $d.Range(1612, 1716)
$fpMain.Text = "We will not have fpStudents for two weeks because of Christmas and New Years Day – we will be back on Wednesday, Jan. 8"

# ---------------------------------------------------------------------------
# Fusion promo paragraph: rewrite with the "RIGHT NOW" / "LAST CHANCE" pitch.
# Using the original-document offsets is safe because this paragraph comes
# before the fpStudents paragraph edited above.
# ---------------------------------------------------------------------------

# Tail: "to register!" -> "to" + "night for their LAST CHANCE ... midnight!"
# (only touch the part after "to", i.e. " register!")
$fusionTail = $d.Range(1500, 1510)
$fusionTail.Text = "night for their LAST CHANCE to lock in the super-low `$49 price. The cost will increase tonight at midnight!"

# Insert "Fusion" right before " – they need to go to "
$fusionInsertPoint = $d.Range(1462, 1462)
$fusionInsertPoint.InsertBefore("Fusion")

# Replace "Fusion is next month" with "RIGHT NOW is the best time to register for "
$fusionMiddle = $d.Range(1442, 1462)
$fusionMiddle.Text = "RIGHT NOW is the best time to register for "

# ---------------------------------------------------------------------------
# Move the _GoBack bookmark from the end of the fpStudents paragraph to the
# end of the (now rewritten) Fusion paragraph. _GoBack is a singleton
# bookmark, so re-adding it at the new location removes it from the old one.
# ---------------------------------------------------------------------------

$fusionPara = $d.Paragraphs(25)
$fusionRange = $fusionPara.Range.Duplicate
$fusionRange.MoveEnd(1, -1)
$bmEnd = $fusionRange.End
$bmRange = $d.Range($bmEnd - 1, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
